$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")

# --- Formatting first (copy existing cell formats into the cells that will
# hold the new rows, before the contents are overwritten) ---

# The old "Summa" row (25) moves down to row 28: carry its bold formatting.
$ws.Range("A25").Copy()
$ws.Range("A28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("B25").Copy()
$ws.Range("B28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# B25 no longer holds the bold total formula - drop its old bold formatting
# so it reverts to the plain/default style used by the rest of column B.
$ws.Range("B25").ClearFormats()

# The three freshly inserted data rows (25-27) get the same date format
# used by column A throughout the rest of the table.
$ws.Range("A24").Copy()
$ws.Range("A25:A27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Now write the new values/formulas ---

$ws.Range("A25").Value = 45320
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "Uppdatera uppgift + Test"

$ws.Range("A26").Value = 45320
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Radera uppgift + Test"

$ws.Range("A27").Value = 45322
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Hämta tidsvy + Test"

# The "Summa" totals row, now on row 28
$ws.Range("A28").Value = "Summa"
$ws.Range("B28").Formula = "=SUBTOTAL(109,B2:B26)"

# --- Expand the Excel table to cover the new rows ---
$lo = $ws.ListObjects.Item("Tabell2")
$lo.Resize($ws.Range("A1:C28"))

# --- Update the selected cell shown in the sheet view ---
$ws.Range("E24").Select() | Out-Null
